$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.081.43'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.27%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.910.22'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.86'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4605'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3873'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07827'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9887'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.92'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.914.33'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.760'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.007'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07045'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.55'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009928'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.04'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.084.59'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.338'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.12'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.122.55'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.083'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.86%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.07'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.38'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.885'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '118.54'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.856'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09310'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8829'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.186'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.315'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.132'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05776'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.169'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02085'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.002'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5689'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.643'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1807'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.720'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.88'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.70%  '

$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002809'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +86.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5325'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.188'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06938'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.838'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.555'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.50'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.60%  '

